$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new quarterly data row (row 78) after the last existing row (77).
# Column A holds a text-like date label ("01-07-2021"). Entering it directly
# via .Value causes Excel to auto-convert it into a date serial number, so we
# write it as a formula that evaluates to the literal text, then paste the
# result back as a plain value (this preserves it as text instead of a date).
$ws.Range("A78").Formula = "=""01-07-2021"""
$ws.Range("A78").Copy()
$ws.Range("A78").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("C78").Value = 3.01
$ws.Range("D78").Value = 4.06
$ws.Range("E78").Value = 4.78
